# Update the date/time format of the "Execution Time" column values
# from "MM/DD/YYYY hh:mm:ss AM/PM" to "DD/MM/YYYY hh:mm:ss AM/PM"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = "28/03/2025 10:42:39 AM"
$ws.Range("E3").Value = "28/03/2025 10:42:43 AM"
$ws.Range("E4").Value = "28/03/2025 10:42:45 AM"
$ws.Range("E5").Value = "28/03/2025 10:42:47 AM"
$ws.Range("E6").Value = "28/03/2025 10:42:47 AM"
$ws.Range("E7").Value = "28/03/2025 10:42:47 AM"
$ws.Range("E8").Value = "28/03/2025 10:42:47 AM"
